$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing Text type (so numeric-looking
# strings like "144.80" or "1.00" keep their exact textual formatting,
# matching the original inline-string cells) without permanently altering
# the cell style (captures & restores the existing style afterwards).
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '65.274.60'
$ws.Range("E2").Value = '  +5.03%  '

Set-TextValue $ws.Range("D3") '3.112.51'
$ws.Range("E3").Value = '  +3.12%  '

$ws.Range("E4").Value = '  +0.10%  '

Set-TextValue $ws.Range("D5") '563.27'
$ws.Range("E5").Value = '  +3.86%  '

Set-TextValue $ws.Range("D6") '144.80'
$ws.Range("E6").Value = '  +8.79%  '

$ws.Range("E7").Value = '  -0.07%  '

Set-TextValue $ws.Range("D8") '3.111.17'
$ws.Range("E8").Value = '  +3.24%  '

Set-TextValue $ws.Range("D9") '0.502'
$ws.Range("E9").Value = '  +2.22%  '

$ws.Range("B10").Value = 'Toncoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D10") '6.47'
$ws.Range("E10").Value = '  +5.82%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range("D11") '0.153'
$ws.Range("E11").Value = '  +4.14%  '

Set-TextValue $ws.Range("D12") '0.471'
$ws.Range("E12").Value = '  +5.61%  '

Set-TextValue $ws.Range("D13") '0.0000232'
$ws.Range("E13").Value = '  +4.93%  '

Set-TextValue $ws.Range("D14") '35.55'
$ws.Range("E14").Value = '  +3.54%  '

Set-TextValue $ws.Range("D15") '3.613.73'
$ws.Range("E15").Value = '  +3.36%  '

Set-TextValue $ws.Range("D16") '65.269.92'
$ws.Range("E16").Value = '  +5.24%  '

Set-TextValue $ws.Range("D17") '3.113.02'
$ws.Range("E17").Value = '  +3.52%  '

$ws.Range("E18").Value = '  +0.91%  '

$ws.Range("E19").Value = '  +2.49%  '

Set-TextValue $ws.Range("D20") '483.58'
$ws.Range("E20").Value = '  +0.71%  '

Set-TextValue $ws.Range("D21") '13.86'
$ws.Range("E21").Value = '  +4.77%  '

Set-TextValue $ws.Range("D22") '0.683'
$ws.Range("E22").Value = '  +1.67%  '

Set-TextValue $ws.Range("D23") '7.60'
$ws.Range("E23").Value = '  +7.90%  '

Set-TextValue $ws.Range("D24") '13.58'
$ws.Range("E24").Value = '  +12.14%  '

Set-TextValue $ws.Range("D25") '81.54'
$ws.Range("E25").Value = '  +0.96%  '

Set-TextValue $ws.Range("D26") '1.00'
$ws.Range("E26").Value = '  -0.11%  '

Set-TextValue $ws.Range("D27") '2.81'
$ws.Range("E27").Value = '  +3.74%  '

Set-TextValue $ws.Range("D28") '8.21'
$ws.Range("E28").Value = '  +6.15%  '

$ws.Range("E29").Value = '  +7.32%  '

Set-TextValue $ws.Range("D30") '0.999'
$ws.Range("E30").Value = '  +0.21%  '

Set-TextValue $ws.Range("D31") '26.33'
$ws.Range("E31").Value = '  +2.46%  '

$ws.Range("E32").Value = '  +3.27%  '

$ws.Range("E33").Value = '  +6.35%  '

Set-TextValue $ws.Range("D34") '5.68'
$ws.Range("E34").Value = '  +0.41%  '

Set-TextValue $ws.Range("D35") '6.21'
$ws.Range("E35").Value = '  +6.01%  '

Set-TextValue $ws.Range("D36") '55.25'
$ws.Range("E36").Value = '  +0.55%  '

Set-TextValue $ws.Range("D37") '473.07'
$ws.Range("E37").Value = '  +4.55%  '

Set-TextValue $ws.Range("D38") '0.0414'
$ws.Range("E38").Value = '  +7.77%  '

Set-TextValue $ws.Range("D39") '0.0838'
$ws.Range("E39").Value = '  +4.90%  '

Set-TextValue $ws.Range("D40") '2.94'
$ws.Range("E40").Value = '  +20.70%  '

Set-TextValue $ws.Range("D41") '3.009.38'
$ws.Range("E41").Value = '  -4.66%  '

$ws.Range("B42").Value = 'Cosmos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D42") '8.29'
$ws.Range("E42").Value = '  +2.61%  '

$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D43") '0.116'
$ws.Range("E43").Value = '  -1.39%  '

Set-TextValue $ws.Range("D44") '28.24'
$ws.Range("E44").Value = '  +6.82%  '

Set-TextValue $ws.Range("D45") '0.262'
$ws.Range("E45").Value = '  +7.58%  '

Set-TextValue $ws.Range("D46") '2.18'
$ws.Range("E46").Value = '  +11.23%  '

$ws.Range("E47").Value = '  +0.03%  '

Set-TextValue $ws.Range("D48") '0.113'
$ws.Range("E48").Value = '  +3.47%  '

Set-TextValue $ws.Range("D49") '0.0₃0539'
$ws.Range("E49").Value = '  +8.62%  '

Set-TextValue $ws.Range("D50") '115.65'
$ws.Range("E50").Value = '  +1.22%  '

Set-TextValue $ws.Range("D51") '2.09'
$ws.Range("E51").Value = '  +3.18%  '

